$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.899.18'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.113.53'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.49'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.31'
$ws.Range("E6").Value = '  +1.04%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.110.87'
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.49'
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.483'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.15'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.629.56'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.882.93'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("E18").Value = '  -1.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.114.15'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.42'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '475.83'
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("E22").Value = '  -0.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.91'
$ws.Range("E23").Value = '  +4.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.50'
$ws.Range("E24").Value = '  +4.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.97'
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("E26").Value = '  -2.62%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.91'
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("E31").Value = '  -0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.53'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0937'
$ws.Range("E34").Value = '  -7.85%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.976'
$ws.Range("E37").Value = '  -3.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.40'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.07'
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.310'
$ws.Range("E41").Value = '  -2.95%  '
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.803.20'
$ws.Range("E45").Value = '  -2.24%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.59'
$ws.Range("E46").Value = '  -11.56%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '380.83'
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.01'
$ws.Range("E48").Value = '  +0.65%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.66'
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("E51").Value = '  -2.34%  '
